# Apply cryptos list update (Wed Jan 17 20:53:00 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.828.81"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.540.50"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'309.81"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "'36.42"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "'7.37"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "2.927.28"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'15.84"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").Value = "2.504.24"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "42.797.11"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'6.79"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "'12.32"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "'69.37"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'244.66"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").Value = "'2.91"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'26.19"
$ws.Range("E27").Value = "  -4.01%  "
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("D29").Value = "'39.57"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "'10.23"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "'5.81"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "'156.03"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0795"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.63"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'2.04"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'18.43"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.20"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.113"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.119"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'4.32"
$ws.Range("E40").Value = "  +9.32%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.17"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0300"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.977.23"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.92"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.782.65"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.194"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").Value = "'0.864"
$ws.Range("E49").Value = "  +10.41%  "
$ws.Range("D50").Value = "'80.97"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'72.87"
$ws.Range("E51").Value = "  -3.07%  "
